# Auto-generated edit script: applies the row-swap changes described by the diff.
# Column A (Sl No) is left untouched; columns B-G are updated to their new values
# (adjacent-row swaps / 3-row rotation), and M1 (Run Date serial) is bumped by 1 day.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value2 = 45961

$ws.Range("B313").Value2 = 57854

$ws.Range("B314").Value2 = 62997

$ws.Range("B316").Value2 = 63565
$ws.Range("E316").Value2 = 109.19
$ws.Range("F316").Value2 = 60
$ws.Range("G316").Value2 = 6162.6

$ws.Range("B317").Value2 = 61610
$ws.Range("D317").Value2 = 102.71
$ws.Range("E317").Value2 = 122.71
$ws.Range("F317").Value2 = -58
$ws.Range("G317").Value2 = -5957.18

$ws.Range("B318").Value2 = 57077
$ws.Range("D318").Value2 = 93.08
$ws.Range("E318").Value2 = 111.2
$ws.Range("F318").Value2 = 0
$ws.Range("G318").Value2 = 0

$ws.Range("B346").Value2 = 63520
$ws.Range("E346").Value2 = 153.4
$ws.Range("F346").Value2 = 85
$ws.Range("G346").Value2 = 12263.8

$ws.Range("B347").Value2 = 55373
$ws.Range("E347").Value2 = 163.62
$ws.Range("F347").Value2 = -94
$ws.Range("G347").Value2 = -13562.32

$ws.Range("B350").Value2 = 63571
$ws.Range("F350").Value2 = 12
$ws.Range("G350").Value2 = 1721.76

$ws.Range("B352").Value2 = 63531
$ws.Range("F352").Value2 = 80
$ws.Range("G352").Value2 = 11478.4

$ws.Range("B372").Value2 = 57885
$ws.Range("E372").Value2 = 62.28
$ws.Range("F372").Value2 = 0
$ws.Range("G372").Value2 = 0

$ws.Range("B373").Value2 = 63652
$ws.Range("E373").Value2 = 55.42
$ws.Range("F373").Value2 = 159
$ws.Range("G373").Value2 = 8288.67

$ws.Range("B379").Value2 = 65514
$ws.Range("F379").Value2 = 0
$ws.Range("G379").Value2 = 0

$ws.Range("B380").Value2 = 63564
$ws.Range("F380").Value2 = 27
$ws.Range("G380").Value2 = 3483.27

$ws.Range("B382").Value2 = 63560
$ws.Range("E382").Value2 = 134.87
$ws.Range("F382").Value2 = 1
$ws.Range("G382").Value2 = 126.86

$ws.Range("B383").Value2 = 60325
$ws.Range("E383").Value2 = 151.57
$ws.Range("F383").Value2 = -102
$ws.Range("G383").Value2 = -12939.72

$ws.Range("B389").Value2 = 57817

$ws.Range("B390").Value2 = 62865

$ws.Range("B400").Value2 = 57835
$ws.Range("F400").Value2 = 0
$ws.Range("G400").Value2 = 0

$ws.Range("B401").Value2 = 62933
$ws.Range("F401").Value2 = 100
$ws.Range("G401").Value2 = 5913

$ws.Range("B419").Value2 = 57856
$ws.Range("F419").Value2 = 0
$ws.Range("G419").Value2 = 0

$ws.Range("B420").Value2 = 63007
$ws.Range("F420").Value2 = 822
$ws.Range("G420").Value2 = 140833.26

$ws.Range("B431").Value2 = 63102
$ws.Range("C431").Value2 = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F431").Value2 = 4
$ws.Range("G431").Value2 = 237.88

$ws.Range("B432").Value2 = 53082
$ws.Range("C432").Value2 = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F432").Value2 = 0
$ws.Range("G432").Value2 = 0

$ws.Range("B536").Value2 = 47097
$ws.Range("D536").Value2 = 112.28
$ws.Range("E536").Value2 = 134.16
$ws.Range("F536").Value2 = 15
$ws.Range("G536").Value2 = 1684.2

$ws.Range("B537").Value2 = 58047
$ws.Range("D537").Value2 = 105.54
$ws.Range("E537").Value2 = 126.1
$ws.Range("F537").Value2 = 43
$ws.Range("G537").Value2 = 4538.22

$ws.Range("B583").Value2 = 53263
$ws.Range("E583").Value2 = 15.29
$ws.Range("F583").Value2 = -309
$ws.Range("G583").Value2 = -3958.29

$ws.Range("B584").Value2 = 65066
$ws.Range("E584").Value2 = 13.61
$ws.Range("F584").Value2 = 221
$ws.Range("G584").Value2 = 2831.01

$ws.Range("B586").Value2 = 45695
$ws.Range("E586").Value2 = 23.58
$ws.Range("F586").Value2 = -36
$ws.Range("G586").Value2 = -710.28

$ws.Range("B587").Value2 = 64915
$ws.Range("E587").Value2 = 20.98
$ws.Range("F587").Value2 = 2
$ws.Range("G587").Value2 = 39.46

$ws.Range("B593").Value2 = 45718
$ws.Range("E593").Value2 = 19.38
$ws.Range("F593").Value2 = -294
$ws.Range("G593").Value2 = -4768.68

$ws.Range("B594").Value2 = 64927
$ws.Range("E594").Value2 = 17.26
$ws.Range("F594").Value2 = 264
$ws.Range("G594").Value2 = 4282.08

$ws.Range("B601").Value2 = 64919
$ws.Range("E601").Value2 = 27.97
$ws.Range("F601").Value2 = 172
$ws.Range("G601").Value2 = 4523.6

$ws.Range("B602").Value2 = 45702
$ws.Range("E602").Value2 = 31.43
$ws.Range("F602").Value2 = -215
$ws.Range("G602").Value2 = -5654.5

$ws.Range("B720").Value2 = 64830
$ws.Range("E720").Value2 = 34.9
$ws.Range("F720").Value2 = 114
$ws.Range("G720").Value2 = 3742.62

$ws.Range("B721").Value2 = 60022
$ws.Range("E721").Value2 = 37.22
$ws.Range("F721").Value2 = -113
$ws.Range("G721").Value2 = -3709.79

$ws.Range("B859").Value2 = 63150
$ws.Range("D859").Value2 = 75.68000000000001
$ws.Range("E859").Value2 = 80.45
$ws.Range("F859").Value2 = 142
$ws.Range("G859").Value2 = 10746.56

$ws.Range("B860").Value2 = 61428
$ws.Range("D860").Value2 = 69.16
$ws.Range("E860").Value2 = 73.52
$ws.Range("F860").Value2 = 1
$ws.Range("G860").Value2 = 69.16

